$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet ---
$ws.Name = "Login"

# --- Cell values (ordered to reproduce the original sharedStrings index order) ---
$ws.Range("A8").Value = "ralph@gmail.com"
$ws.Range("A2").Value = "ralph@gmail.com1"
$ws.Range("A3").Value = "ralph@gmail.com2"
$ws.Range("A7").Value = "ralph@gmail.com7"
$ws.Range("B4").Value = "Pass123$"
$ws.Range("A5").Value = "     "
$ws.Range("A6").Value = "ralphgmail.com6"
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"

$ws.Range("B2").Value = 1234567
$ws.Range("B5").Value = "     "
$ws.Range("B6").Value = "Pass123$"
$ws.Range("B7").Value = 123
$ws.Range("B8").Value = "Pass123$"
# B3 stays blank but still picks up the bordered/centered number style below.

# --- Whole table formatting: thin box border + centered values ---
$all = $ws.Range("A1:B8")
$all.Borders.LineStyle = 1
$all.HorizontalAlignment = -4108

# --- Header row formatting: bold font + yellow fill (border/alignment already applied) ---
$header = $ws.Range("A1:B1")
$header.Font.Bold = $true
$header.Interior.Color = 65535

# --- Column widths (bestFit-style) ---
$ws.Columns.Item(1).ColumnWidth = 17
$ws.Columns.Item(2).ColumnWidth = 8.6

# --- Selection / view ---
$ws.Range("A1:B8").Select() | Out-Null

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
